$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87

# Columns A-D hold text that Excel would otherwise auto-convert
# (dates/numbers). Force text storage, then strip the number-format
# override so the new row matches the plain (unstyled) look of the
# existing data rows.
function Set-TextValue($col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue 1 "2023-06-28"
Set-TextValue 2 "21:41:15"
Set-TextValue 3 "Wednesday"
Set-TextValue 4 "26"

$ws.Cells.Item($row, 5).Value = 123106
$ws.Cells.Item($row, 6).Value = 134304
$ws.Cells.Item($row, 7).Value = 163961
$ws.Cells.Item($row, 8).Value = 134133
$ws.Cells.Item($row, 9).Value = 177220
$ws.Cells.Item($row, 10).Value = 115012
$ws.Cells.Item($row, 11).Value = 204384
$ws.Cells.Item($row, 12).Value = 226522
$ws.Cells.Item($row, 13).Value = 176244
$ws.Cells.Item($row, 14).Value = 104450
$ws.Cells.Item($row, 15).Value = 39773
$ws.Cells.Item($row, 16).Value = 33728
$ws.Cells.Item($row, 17).Value = 52451
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36141
$ws.Cells.Item($row, 20).Value = -1
